# The edit re-orders the species-record rows 4-15 on the "Artfynd" sheet:
# each destination row ends up holding the *entire* row of data (every
# populated column A:AY) that used to live at a different source row.
# Because the permutation contains cycles, every source row is first
# snapshotted into a scratch row far below the used range, then each
# snapshot is written back into its final destination, and finally the
# scratch rows are wiped.
#
# A plain "A:AY" range Copy forces Excel to materialise a (blank) cell in
# every column of the destination, even columns that had no cell at all
# in the source (e.g. J, M:O, X, AC, AF, AH:AS, AU:AV are never used in
# these rows, and column L only has a - always empty - cell on a couple
# of rows). So after the bulk copy we explicitly clear those columns back
# out wherever they must stay/become entirely absent, to reproduce the
# exact sparse layout the source file has.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (row numbers refer to the original,
# pre-edit layout)
$mapping = @{
    4  = 13
    5  = 4
    6  = 14
    7  = 5
    8  = 6
    9  = 7
    10 = 8
    11 = 9
    12 = 15
    13 = 10
    14 = 11
    15 = 12
}

# Source rows whose "L" (Kön) column actually holds a (blank) cell before
# the edit. Every other row has no L cell at all.
$rowsWithL = @(7, 9)

$firstCol = "A"
$lastCol  = "AY"
$scratchOffset = 100   # scratch rows live at sourceRow + 100 (104-115)

# columns that are completely unused (no cell at all) on every one of
# these rows, grouped into contiguous ranges
$alwaysEmptyRanges = @("J:J", "M:O", "X:X", "AC:AC", "AF:AF", "AH:AS", "AU:AV")

# 1) snapshot every involved source row into its scratch row
foreach ($srcRow in ($mapping.Values | Sort-Object -Unique)) {
    $srcRange = $ws.Range("$firstCol$srcRow`:$lastCol$srcRow")
    $bufRow = $srcRow + $scratchOffset
    $bufRange = $ws.Range("$firstCol$bufRow`:$lastCol$bufRow")
    $srcRange.Copy($bufRange)
}

# 2) write each destination row from its scratch snapshot, then restore
#    the correct sparse layout for that row
foreach ($destRow in ($mapping.Keys | Sort-Object)) {
    $srcRow = $mapping[$destRow]
    $bufRow = $srcRow + $scratchOffset
    $bufRange = $ws.Range("$firstCol$bufRow`:$lastCol$bufRow")
    $destRange = $ws.Range("$firstCol$destRow`:$lastCol$destRow")
    $bufRange.Copy($destRange)

    foreach ($colRange in $alwaysEmptyRanges) {
        $parts = $colRange.Split(":")
        $ws.Range("$($parts[0])$destRow`:$($parts[1])$destRow").ClearContents()
    }

    if ($rowsWithL -notcontains $srcRow) {
        $ws.Range("L$destRow").ClearContents()
    }
}

# 3) clean up the scratch rows so the sheet's used range is unaffected
foreach ($srcRow in ($mapping.Values | Sort-Object -Unique)) {
    $bufRow = $srcRow + $scratchOffset
    $bufRange = $ws.Range("$firstCol$bufRow`:$lastCol$bufRow")
    $bufRange.Clear()
}
